$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: Aris, 2024-02-20 (45342), 8:00 -> 9:35 ---
$ws.Range("A26").Value = "Aris"

$ws.Range("B25").Copy() | Out-Null
$ws.Range("B26").PasteSpecial(-4122) | Out-Null
$ws.Range("B26").Value = 45342

$ws.Range("C26").Formula = "=8"
$ws.Range("D26").Formula = "=9+35/60"

# --- Row 27: Viki, 2024-02-20 (45342), 8:00 -> 9:35, What?: Sprites ---
$ws.Range("A27").Value = "Viki"

$ws.Range("B25").Copy() | Out-Null
$ws.Range("B27").PasteSpecial(-4122) | Out-Null
$ws.Range("B27").Value = 45342

$ws.Range("C27").Formula = "=8"
$ws.Range("D27").Formula = "=9+35/60"
$ws.Range("F27").Value = "Sprites"

# Duration column ended up auto-sized after the new entries were typed in.
$ws.Columns("E").ColumnWidth = 7.83

# Selection ends up resting on F28 after entering the new row of data.
$ws.Range("F28").Select() | Out-Null
